# Weekly fruit/vegetable price update: insert two new daily records at the
# top of the historical log (rows 463-464), shifting all existing records
# below down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 463; everything from 463 down shifts to 465+.
$ws.Rows("463:464").Insert()

# New row 463
$ws.Cells.Item(463, 1).Value = 11
$ws.Cells.Item(463, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(463, 3).Value = "Bíobío"
$ws.Cells.Item(463, 4).Value = 45275
$ws.Cells.Item(463, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(463, 5).Value = 8
$ws.Cells.Item(463, 6).Value = 100114001
$ws.Cells.Item(463, 7).Value = "Papa"
$ws.Cells.Item(463, 8).Value = "Asterix"
$ws.Cells.Item(463, 9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(463, 10).Value = 150
$ws.Cells.Item(463, 11).Value = 23000
$ws.Cells.Item(463, 12).Value = 23000
$ws.Cells.Item(463, 13).Value = 23000
$ws.Cells.Item(463, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(463, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(463, 16).Value = 920
$ws.Cells.Item(463, 17).Value = 25
$ws.Cells.Item(463, 18).Value = "Hortaliza"

# New row 464
$ws.Cells.Item(464, 1).Value = 11
$ws.Cells.Item(464, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(464, 3).Value = "Bíobío"
$ws.Cells.Item(464, 4).Value = 45275
$ws.Cells.Item(464, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(464, 5).Value = 8
$ws.Cells.Item(464, 6).Value = 100114001
$ws.Cells.Item(464, 7).Value = "Papa"
$ws.Cells.Item(464, 8).Value = "Asterix"
$ws.Cells.Item(464, 9).Value = "1a nueva(o)"
$ws.Cells.Item(464, 10).Value = 150
$ws.Cells.Item(464, 11).Value = 20000
$ws.Cells.Item(464, 12).Value = 20000
$ws.Cells.Item(464, 13).Value = 20000
$ws.Cells.Item(464, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(464, 15).Value = "Región del Maule"
$ws.Cells.Item(464, 16).Value = 800
$ws.Cells.Item(464, 17).Value = 25
$ws.Cells.Item(464, 18).Value = "Hortaliza"
